# Applies the "StructureDefinition-covered-text" update:
#  - rebrand from ibm.com/Alvearie to linuxforhealth.org/LinuxForHealth
#  - bump version 7.0.0 -> 8.0.0
#  - bump the publication date
#  - clear the stray "Constraint(s)" text that had leaked onto the
#    top-level Extension row in the Elements table (it only belongs on
#    the Extension.extension row)

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/covered-text"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the top-level "Extension" element; its "Constraint(s)" column
# (AI) incorrectly carried the ele-1/ext-1 constraint text that belongs
# on the "Extension.extension" row (row 4) instead. Clear it.
$elements.Range("AI2").ClearContents()

# Row 5 is the "Extension.url" element; per the FHIR convention its Fixed
# Value (column Q) mirrors the StructureDefinition's own canonical URL, so
# it needs the same rebrand as the Metadata sheet's URL above.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/covered-text"
